$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds the "last changed" date serial for each record.
# The update bumps this date by one day (45177 -> 45178) for every data row
# (rows 2 through 270).
$ws.Range("C2:C270").Value2 = 45178
